$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 78: 2025-XX-XX (serial 45939), site 四方坪站 (shared string idx 2)
$ws.Range("A78").Value = 45939
$ws.Range("B78").Value = "四方坪站"
$ws.Range("C78").Formula = "=18835/126"
$ws.Range("D78").Formula = "=C78/(24*60)"
$ws.Range("E78").Formula = "=10987.23/126"
$ws.Range("F78").Formula = "=3818.52/126"
$ws.Range("G78").Formula = "=10987.23/(18835/60)"
$ws.Range("H78").Formula = "=460/126"

# Row 79: serial 45939, site 高岭站 (shared string idx 3)
$ws.Range("A79").Value = 45939
$ws.Range("B79").Value = "高岭站"
$ws.Range("C79").Formula = "=7356/36"
$ws.Range("D79").Formula = "=C79/(24*60)"
$ws.Range("E79").Formula = "=5141.28/36"
$ws.Range("F79").Formula = "=1326.75/36"
$ws.Range("G79").Formula = "=5141.28/(7356/60)"
$ws.Range("H79").Formula = "=179/36"

# Row 80: serial 45940, site 四方坪站
$ws.Range("A80").Value = 45940
$ws.Range("B80").Value = "四方坪站"
$ws.Range("C80").Formula = "=20047/126"
$ws.Range("D80").Formula = "=C80/(24*60)"
$ws.Range("E80").Formula = "=10256.78/126"
$ws.Range("F80").Formula = "=3626.33/126"
$ws.Range("G80").Formula = "=10256.78/(20047/60)"
$ws.Range("H80").Formula = "=440/126"

# Row 81: serial 45940, site 高岭站
$ws.Range("A81").Value = 45940
$ws.Range("B81").Value = "高岭站"
$ws.Range("C81").Formula = "=7592/36"
$ws.Range("D81").Formula = "=C81/(24*60)"
$ws.Range("E81").Formula = "=5154.54/36"
$ws.Range("F81").Formula = "=1307.95/36"
$ws.Range("G81").Formula = "=5154.54/(7592/60)"
$ws.Range("H81").Formula = "=191/36"

# Match the view state captured in the saved workbook
$ws.Range("H82").Select()
$excel.ActiveWindow.ScrollRow = 67
